$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).NumberFormat = "General"
    $ws.Range($ref).Style = "Normal"
}

$ws.Range('D2').Value = '52.311.24'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '3.014.34'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  -0.11%  '
Set-TextValue 'D5' '355.27'
$ws.Range('E5').Value = '  +0.98%  '
Set-TextValue 'D6' '108.27'
$ws.Range('E6').Value = '  -2.94%  '
Set-TextValue 'D7' '0.563'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -1.73%  '
Set-TextValue 'D10' '38.57'
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('E11').Value = '  +2.06%  '
Set-TextValue 'D12' '0.0862'
$ws.Range('E12').Value = '  -3.62%  '
Set-TextValue 'D13' '19.31'
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('D14').Value = '3.485.14'
$ws.Range('E14').Value = '  +2.11%  '
Set-TextValue 'D15' '7.70'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').Value = '3.001.69'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('E17').Value = '  +3.11%  '
$ws.Range('D18').Value = '52.329.24'
$ws.Range('E18').Value = '  +0.80%  '
Set-TextValue 'D19' '3.54'
$ws.Range('E19').Value = '  +8.18%  '
Set-TextValue 'D20' '7.55'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('E21').Value = '  -5.19%  '
$ws.Range('E22').Value = '  -1.13%  '
Set-TextValue 'D23' '69.68'
$ws.Range('E23').Value = '  -2.38%  '
Set-TextValue 'D24' '265.69'
$ws.Range('E24').Value = '  -2.69%  '
Set-TextValue 'D25' '2.75'
$ws.Range('E25').Value = '  -0.95%  '
Set-TextValue 'D26' '0.179'
$ws.Range('E26').Value = '  -1.93%  '
Set-TextValue 'D27' '7.75'
$ws.Range('E27').Value = '  +4.70%  '
Set-TextValue 'D28' '27.05'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('E31').Value = '  +2.46%  '
Set-TextValue 'D32' '10.35'
$ws.Range('E32').Value = '  -3.71%  '
Set-TextValue 'D33' '36.47'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('E34').Value = '  +15.15%  '
Set-TextValue 'D35' '50.99'
$ws.Range('E35').Value = '  -3.73%  '
Set-TextValue 'D36' '0.0445'
$ws.Range('E36').Value = '  -0.74%  '
Set-TextValue 'D37' '0.998'
$ws.Range('E37').Value = '  -0.08%  '
Set-TextValue 'D38' '3.24'
$ws.Range('E38').Value = '  -4.43%  '
$ws.Range('E39').Value = '  -2.08%  '
Set-TextValue 'D40' '17.97'
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  +2.44%  '
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('E43').Value = '  -2.76%  '
Set-TextValue 'D44' '123.57'
$ws.Range('E44').Value = '  +9.06%  '
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('D46').Value = '2.129.93'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('E47').Value = '  -3.89%  '
Set-TextValue 'D48' '2.39'
$ws.Range('E48').Value = '  -5.66%  '
$ws.Range('D49').Value = '3.311.31'
$ws.Range('E49').Value = '  +2.15%  '
Set-TextValue 'D50' '0.249'
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('E51').Value = '  -1.41%  '
